$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H75").Value = 29623
$ws.Range("J75").Value = 29623
$ws.Range("L75").Value = 29623
$ws.Range("N75").Value = -31495
$ws.Range("H78").Value = 29623
$ws.Range("J78").Value = 29623
$ws.Range("L78").Value = 88869
$ws.Range("N78").Value = -98229
$ws.Range("H106").Value = 1067.25
$ws.Range("I106").Value = 600
$ws.Range("J106").Value = 1223
$ws.Range("K106").Value = 600
$ws.Range("L106").Value = 1223
$ws.Range("M106").Value = 31
$ws.Range("N106").Value = -2485
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200
$ws.Range("H138").Value = 4139.467
$ws.Range("I138").Value = 2161.5
$ws.Range("K138").Value = 6484.5
$ws.Range("M138").Value = -1344.5

$ws = $wb.Worksheets("ARM")
$ws.Range("H34").Value = 28705.6
$ws.Range("I34").Value = 10500
$ws.Range("J34").Value = 33257
$ws.Range("K34").Value = 10500
$ws.Range("L34").Value = 33257
$ws.Range("M34").Value = -10229
$ws.Range("N34").Value = -33799
$ws.Range("H61").Value = 3649.484
$ws.Range("I61").Value = 734.1539
$ws.Range("J61").Value = 5755
$ws.Range("K61").Value = 734.1539
$ws.Range("L61").Value = 5755
$ws.Range("M61").Value = -522.1539
$ws.Range("N61").Value = -6179
$ws.Range("H92").Value = 31975.8
$ws.Range("J92").Value = 31975.8
$ws.Range("L92").Value = 31975.8
$ws.Range("N92").Value = -36967.8
$ws.Range("H101").Value = 38881
$ws.Range("J101").Value = 38881
$ws.Range("L101").Value = 38881
$ws.Range("N101").Value = -45371
$ws.Range("H102").Value = 2094.7576
$ws.Range("I102").Value = 1764.7407
$ws.Range("J102").Value = 3579.8333
$ws.Range("K102").Value = 1764.7407
$ws.Range("L102").Value = 3579.8333
$ws.Range("M102").Value = -142.7407000000001
$ws.Range("N102").Value = -6823.8333
$ws.Range("H119").Value = 38107.2
$ws.Range("I119").Value = 10000
$ws.Range("J119").Value = 45134
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 45134
$ws.Range("M119").Value = -5162
$ws.Range("N119").Value = -54810
$ws.Range("H136").Value = 3649.484
$ws.Range("I136").Value = 734.1539
$ws.Range("J136").Value = 5755
$ws.Range("K136").Value = 2202.4617
$ws.Range("L136").Value = 17265
$ws.Range("M136").Value = 347.5383000000002
$ws.Range("N136").Value = -22365

$ws = $wb.Worksheets("BSM")
$ws.Range("H76").Value = 34600
$ws.Range("I76").Value = 13000
$ws.Range("J76").Value = 40000
$ws.Range("K76").Value = 13000
$ws.Range("L76").Value = 40000
$ws.Range("M76").Value = -12685
$ws.Range("N76").Value = -40630
$ws.Range("H79").Value = 34600
$ws.Range("I79").Value = 13000
$ws.Range("J79").Value = 40000
$ws.Range("K79").Value = 13000
$ws.Range("L79").Value = 40000
$ws.Range("M79").Value = -11908
$ws.Range("N79").Value = -42184
$ws.Range("H94").Value = 2006.3871
$ws.Range("I94").Value = 1873.7142
$ws.Range("J94").Value = 2285
$ws.Range("K94").Value = 1873.7142
$ws.Range("L94").Value = 2285
$ws.Range("M94").Value = -1422.7142
$ws.Range("N94").Value = -3187

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 1857.1428
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -6700
$ws.Range("H41").Value = 36287.5
$ws.Range("J41").Value = 36287.5
$ws.Range("L41").Value = 36287.5
$ws.Range("N41").Value = -37143.5
$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -39492

$ws = $wb.Worksheets("CUL")
$ws.Range("H68").Value = 1935.8049
$ws.Range("I68").Value = 803.5
$ws.Range("J68").Value = 2301.0645
$ws.Range("K68").Value = 2410.5
$ws.Range("L68").Value = 6903.193499999999
$ws.Range("M68").Value = -1599.5
$ws.Range("N68").Value = -8525.193499999999
$ws.Range("H71").Value = 1935.8049
$ws.Range("I71").Value = 803.5
$ws.Range("J71").Value = 2301.0645
$ws.Range("K71").Value = 7231.5
$ws.Range("L71").Value = 20709.5805
$ws.Range("M71").Value = -3175.5
$ws.Range("N71").Value = -28821.5805
$ws.Range("H129").Value = 27214.45
$ws.Range("J129").Value = 37149.855
$ws.Range("L129").Value = 111449.565
$ws.Range("N129").Value = -121449.565
$ws.Range("H130").Value = 2450
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 2500
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 7500
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -17540
$ws.Range("H131").Value = 1718.8462
$ws.Range("I131").Value = 6300
$ws.Range("J131").Value = 1337.0834
$ws.Range("K131").Value = 18900
$ws.Range("L131").Value = 4011.2502
$ws.Range("M131").Value = -13860
$ws.Range("N131").Value = -14091.2502

$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 3716.5
$ws.Range("J80").Value = 4933.3335
$ws.Range("L80").Value = 4933.3335
$ws.Range("N80").Value = -6929.3335
$ws.Range("H83").Value = 3716.5
$ws.Range("J83").Value = 4933.3335
$ws.Range("L83").Value = 24666.6675
$ws.Range("N83").Value = -34650.6675
$ws.Range("H113").Value = 2535.7778
$ws.Range("I113").Value = 1316
$ws.Range("K113").Value = 1316
$ws.Range("M113").Value = 854
$ws.Range("H126").Value = 4369.231
$ws.Range("I126").Value = 2666.6667
$ws.Range("J126").Value = 4880
$ws.Range("K126").Value = 8000.000100000001
$ws.Range("L126").Value = 14640
$ws.Range("M126").Value = -5530.000100000001
$ws.Range("N126").Value = -19580

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 2044.1177
$ws.Range("I7").Value = 1432.8182
$ws.Range("J7").Value = 3164.8333
$ws.Range("K7").Value = 1432.8182
$ws.Range("L7").Value = 3164.8333
$ws.Range("M7").Value = -1320.8182
$ws.Range("N7").Value = -3388.8333
$ws.Range("H16").Value = 90912180
$ws.Range("I16").Value = 111112880
$ws.Range("K16").Value = 111112880
$ws.Range("M16").Value = -111112710
$ws.Range("H22").Value = 1429.1875
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 2074.111
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 2074.111
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -2664.111
$ws.Range("H27").Value = 1429.1875
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 2074.111
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 2074.111
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -2288.111
$ws.Range("H40").Value = 2563.182
$ws.Range("I40").Value = 1496.6666
$ws.Range("J40").Value = 2963.125
$ws.Range("K40").Value = 1496.6666
$ws.Range("L40").Value = 2963.125
$ws.Range("M40").Value = -1360.6666
$ws.Range("N40").Value = -3235.125
$ws.Range("H93").Value = 2537.2
$ws.Range("I93").Value = 1998.4546
$ws.Range("J93").Value = 4018.75
$ws.Range("K93").Value = 1998.4546
$ws.Range("L93").Value = 4018.75
$ws.Range("M93").Value = -750.4546
$ws.Range("N93").Value = -6514.75
$ws.Range("H126").Value = 2044.1177
$ws.Range("I126").Value = 1432.8182
$ws.Range("J126").Value = 3164.8333
$ws.Range("K126").Value = 4298.4546
$ws.Range("L126").Value = 9494.499899999999
$ws.Range("M126").Value = -1828.4546
$ws.Range("N126").Value = -14434.4999
$ws.Range("H132").Value = 3679.8635
$ws.Range("I132").Value = 2750.375
$ws.Range("J132").Value = 4211
$ws.Range("K132").Value = 8251.125
$ws.Range("L132").Value = 12633
$ws.Range("M132").Value = -5721.125
$ws.Range("N132").Value = -17693

$ws = $wb.Worksheets("WVR")
$ws.Range("H69").Value = 25343.666
$ws.Range("J69").Value = 25343.666
$ws.Range("L69").Value = 25343.666
$ws.Range("N69").Value = -26841.666
$ws.Range("H72").Value = 25343.666
$ws.Range("J72").Value = 25343.666
$ws.Range("L72").Value = 76030.99800000001
$ws.Range("N72").Value = -83518.99800000001
$ws.Range("H122").Value = 2926.2632
$ws.Range("I122").Value = 2036
$ws.Range("J122").Value = 5419
$ws.Range("K122").Value = 6108
$ws.Range("L122").Value = 16257
$ws.Range("M122").Value = -3658
$ws.Range("N122").Value = -21157
$ws.Range("H132").Value = 11928.786
$ws.Range("I132").Value = 3875.5
$ws.Range("J132").Value = 22666.5
$ws.Range("K132").Value = 11626.5
$ws.Range("L132").Value = 67999.5
$ws.Range("M132").Value = -9096.5
$ws.Range("N132").Value = -73059.5
